# Add a new submission row ("221118_xgb_external") to the tracker table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel table by one row; this also updates the table ref/autoFilter
# and the sheet dimension automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$newRowIndex = $newRow.Range.Row

# Column A - Date. Set the number format *before* the value so the engine
# reuses the existing built-in date format (numFmtId 14) instead of creating
# a brand-new custom number format.
$ws.Cells.Item($newRowIndex, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item($newRowIndex, 1).Value = "11/18/2022"

# Column B - Name
$ws.Cells.Item($newRowIndex, 2).Value = "221118_xgb_external"

# Column C - Name Ramp
$ws.Cells.Item($newRowIndex, 3).Value = "6_JM_MS"

# Column D - Hand in. Copy the existing "TRUE" text cell so the new cell
# stays a shared-string text value instead of becoming a native boolean.
$ws.Range("D4").Copy()
$ws.Cells.Item($newRowIndex, 4).PasteSpecial(-4163)

# Column E - By
$ws.Cells.Item($newRowIndex, 5).Value = "Maria"

# Move the active selection below the newly added row, like Excel would
# after finishing data entry.
$ws.Cells.Item($newRowIndex + 1, 3).Select()
